$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + volume). Row 41/42: Aptos <-> TheSandbox swap.
# D-column "Price" cells are plain text that can look numeric (e.g. "52.30", "0.9960",
# "24.398.47"); force text format before assignment so Excel keeps the exact string
# (trailing zeros, thousands-dot grouping, etc.) instead of auto-converting to a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.398.47"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.52"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9977"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.56"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9960"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3937"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3920"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.30"
$ws.Range("E9").Value = "  +6.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.406"
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9946"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08585"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.45"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.298"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001334"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.880"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.655.08"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.57"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06943"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.53"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.998"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9957"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.67"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.462.09"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.444"
$ws.Range("E25").Value = "  +3.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.964"
$ws.Range("E26").Value = "  +7.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.53"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.48"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "142.81"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.340"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.140"
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.555"
$ws.Range("E32").Value = "  +3.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.843.44"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.061"
$ws.Range("E34").Value = "  +6.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08237"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.845"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02976"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.12"
$ws.Range("E38").Value = "  +10.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2749"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09264"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7758"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.84"
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.64"
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7106"
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.525"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.136"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9954"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08443"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.452"
$ws.Range("E50").Value = "  +12.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.94"
$ws.Range("E51").Value = "  +1.04%  "
